$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("contacts")

$data = @(
    @("Varsha", "Singhal", "MindTree"),
    @("Nidhi", "Choudhary", "Infosys"),
    @("Upasana", "Sinha", "Cognizant"),
    @("Ruchita", "Kadam", "IBM")
)

$startRow = 5
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

$ws.Range("C8").Select()
